# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PIB 2022 Deflacionado -> PIB 2023 Deflacionado section (rows 2-10)
$ws.Range("A5").Value = "SC"
$ws.Range("A6").Value = "RJ"

$ws.Range("B2").Value = 118174.1116095417
$ws.Range("B3").Value = 76532.28963539573
$ws.Range("B4").Value = 73845.19036585005
$ws.Range("B5").Value = 69959.10090505773
$ws.Range("B6").Value = 67161.88535005336
$ws.Range("B7").Value = 64948.89321994126
$ws.Range("B8").Value = 26006.98661973922
$ws.Range("B9").Value = 51300.70579350938
$ws.Range("B10").Value = 26237.41536180414

$ws.Range("C8").Value = 23

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "PIB 2023 Deflacionado"
}

# Variação (%) 2022/2010 -> Variação (%) 2023/2010 section (rows 11-19)
$ws.Range("A11").Value = "PI"
$ws.Range("A12").Value = "AL"
$ws.Range("A14").Value = "PR"
$ws.Range("A15").Value = "MT"
$ws.Range("A16").Value = "RS"

$ws.Range("B11").Value = 1.537570136346218
$ws.Range("B12").Value = 1.47614098883597
$ws.Range("B13").Value = 1.456636142415122
$ws.Range("B14").Value = 1.361548723808843
$ws.Range("B15").Value = 1.347864888435147
$ws.Range("B16").Value = 1.330588661237085
$ws.Range("B17").Value = 0.91785810919599
$ws.Range("B18").Value = 1.205735709293767
$ws.Range("B19").Value = 1.297733510014661

$ws.Range("C17").Value = 23

for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value = "Variação (%) 2023/2010"
}
